# Data.xlsx edit: "adding the Listner class"
#
# 1. Sheet1!B5 (Date): 18 -> 19
# 2. Sheet2!B3 (Pax 1 Firstname): "Tamil hasan" -> "Tamil"
# 3. Sheet2's active cell/selection moves from B3 to B14
# 4. The active (selected) worksheet changes from Sheet2 to Sheet1,
#    leaving Sheet1's own selection on B5.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws1.Range("B5").Value = 19
$ws2.Range("B3").Value = "Tamil"

# Move Sheet2's selection to B14, then make Sheet1 the active tab with
# its selection on B5 (matching the pre-existing selection there).
$ws2.Range("B14").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("B5").Select() | Out-Null
